# maj template comment à la fin
# The "Comment" column (J) is moved to the end of the primer block (after
# TemplateAmount, i.e. column R). Columns K:R (ForwardPrimerName ..
# TemplateAmount) shift one column to the left, into J:Q.
#
# Before: J=Comment, K=ForwardPrimerName, L=ForwardPrimerSequence,
#         M=ReversePrimerName, N=ReversePrimerSequence, O=ExpectedAmpliconSize,
#         P=PcrProgram, Q=TargetGeneName, R=TemplateAmount
# After:  J=ForwardPrimerName, K=ForwardPrimerSequence, L=ReversePrimerName,
#         M=ReversePrimerSequence, N=ExpectedAmpliconSize, O=PcrProgram,
#         P=TargetGeneName, Q=TemplateAmount, R=Comment

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 5

for ($row = 1; $row -le $lastRow; $row++) {
    # Grab the current values of the Comment .. TemplateAmount block (J:R).
    $comment = $ws.Cells.Item($row, 10).Value2          # J = 10
    $forwardPrimerName = $ws.Cells.Item($row, 11).Value2 # K = 11
    $forwardPrimerSeq  = $ws.Cells.Item($row, 12).Value2 # L = 12
    $reversePrimerName = $ws.Cells.Item($row, 13).Value2 # M = 13
    $reversePrimerSeq  = $ws.Cells.Item($row, 14).Value2 # N = 14
    $expectedAmplicon  = $ws.Cells.Item($row, 15).Value2 # O = 15
    $pcrProgram        = $ws.Cells.Item($row, 16).Value2 # P = 16
    $targetGeneName    = $ws.Cells.Item($row, 17).Value2 # Q = 17
    $templateAmount    = $ws.Cells.Item($row, 18).Value2 # R = 18

    # Write them back shifted: primer block moves left into J:Q, Comment goes to R.
    $ws.Cells.Item($row, 10).Value = $forwardPrimerName
    $ws.Cells.Item($row, 11).Value = $forwardPrimerSeq
    $ws.Cells.Item($row, 12).Value = $reversePrimerName
    $ws.Cells.Item($row, 13).Value = $reversePrimerSeq
    $ws.Cells.Item($row, 14).Value = $expectedAmplicon
    $ws.Cells.Item($row, 15).Value = $pcrProgram
    $ws.Cells.Item($row, 16).Value = $targetGeneName
    $ws.Cells.Item($row, 17).Value = $templateAmount
    $ws.Cells.Item($row, 18).Value = $comment
}
